$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was an unused/blank row above the real table; remove it so the
# existing header + data rows shift up by one.
$ws.Rows.Item(1).Delete()

# Append the fourth example as the new last row.
$ws.Cells.Item(5,1).Value = "Estimador Mensual de Actividad Económica. Números índice, base 2004=100 y variaciones porcentuales."
$ws.Cells.Item(5,2).Value = "Monthly"
$ws.Cells.Item(5,3).Value = "example_cleaning_databases_4"
$ws.Cells.Item(5,4).Value = "https://www.indec.gob.ar/indec/web/Nivel4-Tema-3-9-48"
$ws.Cells.Item(5,5).Value = "Argentina"
$ws.Range("B5:E5").HorizontalAlignment = -4108

$ws.Columns.Item(1).EntireColumn.AutoFit()

[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A3").Select()
